$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = 7
$ws.Range("F5").Value = 0
$ws.Range("F8").Value = 2
